{"js": "// Nationality sample table: correct a handful of figures that had not\n// been updated early enough (bug chase related to nationality).\n//   Great Britain: Frequency 239 -> 238, Percent 54 -> 53, Whole TSTD Frequency 11,796 -> 11,797\n//   France: Frequency 84 -> 85\n//   Total: Whole TSTD Frequency 33,516 -> 33,517\nconst body = context.document.body;\n\nconst replacements = [\n  [\"239\", \"238\"],\n  [\"54\", \"53\"],\n  [\"11,796\", \"11,797\"],\n  [\"84\", \"85\"],\n  [\"33,516\", \"33,517\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchWholeWord: true, matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Nationality sample table: correct a handful of figures that had not\n# been updated early enough (bug chase related to nationality).\n#   Great Britain: Frequency 239 -> 238, Percent 54 -> 53, Whole TSTD Frequency 11,796 -> 11,797\n#   France: Frequency 84 -> 85\n#   Total: Whole TSTD Frequency 33,516 -> 33,517\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# row, col (1-based), old value, new value\n$edits = @(\n    @(3, 2, \"239\", \"238\"),\n    @(3, 3, \"54\", \"53\"),\n    @(3, 4, \"11,796\", \"11,797\"),\n    @(5, 2, \"84\", \"85\"),\n    @(11, 4, \"33,516\", \"33,517\")\n)\n\nforeach ($edit in $edits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $oldValue = $edit[2]\n    $newValue = $edit[3]\n\n    $cell = $t.Cell($row, $col)\n    # Cell.Range.Text carries Word's trailing cell-mark pair (\"\\r\\x07\"); trim\n    # it before comparing against the plain expected value.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -ne $oldValue) {\n        throw \"Unexpected value at row $row, col $col : expected '$oldValue' but found '$current'\"\n    }\n\n    $cell.Range.Text = $newValue\n}\n"}
